$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    115 = 56.6
    183 = 5
    185 = 60.3
    186 = 50.8
    187 = 61.4
    188 = 92.6
    190 = 8.1
    191 = 92.8
    266 = 96.6
    267 = 46.7
    268 = 35.2
    269 = 30.9
    271 = 97.3
    272 = 62.7
    273 = 99
    274 = 44.6
    275 = 63.5
    277 = 63.7
    365 = 33.6
    424 = 74.2
    425 = 97.9
    426 = 88.4
    431 = 130
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 1).Value = $changes[$row]
}
